$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 onto the new
# header cells I1:J1 before setting their values.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(4, 4),
    @(2, 2),
    @(8, 8),
    @(12, 12),
    @(6, 7),
    @(7, 8),
    @(3, 4),
    @(6, 7),
    @(7, 7),
    @(3, 3),
    @(8, 9)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
